# Update imputed values in the RandomForest algorithm result sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E6").Value = 16.35410000000002
$ws.Range("C7").Value = -12.6267
$ws.Range("B9").Value = 6.064699999999993
$ws.Range("C12").Value = -10.791
$ws.Range("C14").Value = -12.6268
$ws.Range("E15").Value = 16.3448
$ws.Range("B18").Value = 6.123899999999999
$ws.Range("B20").Value = 8.911600000000004
$ws.Range("C26").Value = -12.67770000000001
$ws.Range("B27").Value = 6.233200000000002
$ws.Range("C27").Value = -12.1699
$ws.Range("C29").Value = -11.1086
$ws.Range("E33").Value = 17.23760000000001
$ws.Range("B35").Value = 8.7666
$ws.Range("E35").Value = 16.4015
$ws.Range("C37").Value = -12.72840000000001
$ws.Range("C38").Value = -12.3898
$ws.Range("E38").Value = 16.53939999999999
$ws.Range("E43").Value = 17.24030000000001
$ws.Range("E44").Value = 16.80709999999998
$ws.Range("E47").Value = 16.46309999999999
$ws.Range("C51").Value = -11.235
$ws.Range("E51").Value = 17.3232
$ws.Range("C52").Value = -11.22040000000001
$ws.Range("C55").Value = -13.8179
$ws.Range("E57").Value = 16.56280000000002
$ws.Range("E63").Value = 18.75620000000001
$ws.Range("B69").Value = 5.602799999999995
$ws.Range("C69").Value = -11.064
$ws.Range("C70").Value = -11.36
$ws.Range("E70").Value = 17.44790000000001
$ws.Range("B76").Value = 5.773999999999997
$ws.Range("B78").Value = 9.917100000000003
$ws.Range("C81").Value = -12.25319999999999
$ws.Range("B82").Value = 5.077000000000002
$ws.Range("B83").Value = 5.306499999999996
$ws.Range("C83").Value = -13.8132
$ws.Range("E88").Value = 16.4705
$ws.Range("B93").Value = 6.251999999999999
$ws.Range("E99").Value = 16.5807
$ws.Range("C102").Value = -12.72710000000001
